$d = $word.ActiveDocument

# --- Paragraph 1: collapse the three runs ("S" / "dfasdfsdfsdf" /
#     "dsfgsdgfdjgiafkljdfkljas;idei;asdjfk;ljasd;kjads") into one run
#     reading "What is 3+3?" ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$r1.Text = "What is 3+3?"

# --- Paragraph 2 ("F"): becomes "Ans:" ---
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.MoveEnd(1, -1) | Out-Null
$r2.Text = "Ans:"

# --- Paragraph 5 ("Asdflkjsda;kfljasdfkl"): clear its text but keep the
#     bookmark that lives at the end of the paragraph ---
$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$r5.MoveEnd(1, -1) | Out-Null
$r5.Text = ""

# --- Remove paragraph 3 ("Asdf"), including its paragraph mark ---
$p3 = $d.Paragraphs.Item(3)
$d.Range($p3.Range.Start, $p3.Range.End).Delete()

# --- Remove the now-empty paragraph that used to be paragraph 4 ---
$p3b = $d.Paragraphs.Item(3)
$d.Range($p3b.Range.Start, $p3b.Range.End).Delete()

# --- Merge what remains of paragraph 5 (now just the bookmark) up into
#     paragraph 2 by deleting paragraph 2's trailing paragraph mark ---
$p2b = $d.Paragraphs.Item(2)
$d.Range($p2b.Range.End - 1, $p2b.Range.End).Delete()
